{"js": "// Replace the multiplication-problem text in the table cells with the\n// newly generated values. Each old equation string is unique within the\n// document, so a simple search+replace per pair is unambiguous.\nconst replacements = [\n  [\"79\u00d760=4740\", \"38\u00d798=3724\"],\n  [\"22\u00d731=682\", \"79\u00d713=1027\"],\n  [\"94\u00d731=2914\", \"69\u00d774=5106\"],\n  [\"26\u00d742=1092\", \"19\u00d798=1862\"],\n  [\"17\u00d772=1224\", \"87\u00d778=6786\"],\n  [\"51\u00d753=2703\", \"62\u00d770=4340\"],\n  [\"41\u00d761=2501\", \"37\u00d758=2146\"],\n  [\"51\u00d765=3315\", \"23\u00d748=1104\"],\n  [\"33\u00d750=1650\", \"60\u00d742=2520\"],\n  [\"96\u00d736=3456\", \"21\u00d766=1386\"],\n  [\"86\u00d729=2494\", \"49\u00d747=2303\"],\n  [\"13\u00d773=949\", \"56\u00d729=1624\"],\n  [\"48\u00d756=2688\", \"50\u00d766=3300\"],\n  [\"53\u00d764=3392\", \"72\u00d741=2952\"],\n  [\"17\u00d771=1207\", \"46\u00d755=2530\"],\n  [\"99\u00d773=7227\", \"28\u00d776=2128\"],\n  [\"75\u00d741=3075\", \"14\u00d762=868\"],\n  [\"29\u00d726=754\", \"63\u00d753=3339\"],\n  [\"28\u00d763=1764\", \"96\u00d779=7584\"],\n  [\"49\u00d781=3969\", \"15\u00d771=1065\"],\n  [\"36\u00d743=1548\", \"55\u00d759=3245\"],\n  [\"78\u00d779=6162\", \"27\u00d768=1836\"],\n  [\"66\u00d791=6006\", \"74\u00d759=4366\"],\n  [\"12\u00d777=924\", \"34\u00d744=1496\"],\n  [\"12\u00d727=324\", \"33\u00d779=2607\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in the table cells with the\n# newly generated values. Each old equation string is unique within the\n# document, so Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"79\u00d760=4740\"; New = \"38\u00d798=3724\" },\n    @{ Old = \"22\u00d731=682\";  New = \"79\u00d713=1027\" },\n    @{ Old = \"94\u00d731=2914\"; New = \"69\u00d774=5106\" },\n    @{ Old = \"26\u00d742=1092\"; New = \"19\u00d798=1862\" },\n    @{ Old = \"17\u00d772=1224\"; New = \"87\u00d778=6786\" },\n    @{ Old = \"51\u00d753=2703\"; New = \"62\u00d770=4340\" },\n    @{ Old = \"41\u00d761=2501\"; New = \"37\u00d758=2146\" },\n    @{ Old = \"51\u00d765=3315\"; New = \"23\u00d748=1104\" },\n    @{ Old = \"33\u00d750=1650\"; New = \"60\u00d742=2520\" },\n    @{ Old = \"96\u00d736=3456\"; New = \"21\u00d766=1386\" },\n    @{ Old = \"86\u00d729=2494\"; New = \"49\u00d747=2303\" },\n    @{ Old = \"13\u00d773=949\";  New = \"56\u00d729=1624\" },\n    @{ Old = \"48\u00d756=2688\"; New = \"50\u00d766=3300\" },\n    @{ Old = \"53\u00d764=3392\"; New = \"72\u00d741=2952\" },\n    @{ Old = \"17\u00d771=1207\"; New = \"46\u00d755=2530\" },\n    @{ Old = \"99\u00d773=7227\"; New = \"28\u00d776=2128\" },\n    @{ Old = \"75\u00d741=3075\"; New = \"14\u00d762=868\"  },\n    @{ Old = \"29\u00d726=754\";  New = \"63\u00d753=3339\" },\n    @{ Old = \"28\u00d763=1764\"; New = \"96\u00d779=7584\" },\n    @{ Old = \"49\u00d781=3969\"; New = \"15\u00d771=1065\" },\n    @{ Old = \"36\u00d743=1548\"; New = \"55\u00d759=3245\" },\n    @{ Old = \"78\u00d779=6162\"; New = \"27\u00d768=1836\" },\n    @{ Old = \"66\u00d791=6006\"; New = \"74\u00d759=4366\" },\n    @{ Old = \"12\u00d777=924\";  New = \"34\u00d744=1496\" },\n    @{ Old = \"12\u00d727=324\";  New = \"33\u00d779=2607\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
